# Read & Write Using Openpyxl
# Replace the sample 가/나/다/라/마/바 header row on the "Change" sheet with a
# small "people" table (header row + 5 data rows), shifted one column to the
# right so the row index lives in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Change")
$ws.Activate()

# Clear the previous A1:F1 content entirely - new layout only needs A1:D6.
$ws.Cells.Clear()

# Header row (no value in A1 - column A holds the row number starting at row 2)
$ws.Range("B1").Value = "이름"
$ws.Range("C1").Value = "사는 곳"
$ws.Range("D1").Value = "직업"

# Data rows
$data = @(
    @(1, "Elsa", "NewYork", "Lawyer"),
    @(2, "Anna", "LA", "Programmer"),
    @(3, "Olaf", "Washington", "Teacher"),
    @(4, "Kim", "Seoul", "Cook"),
    @(5, "Lee", "Busan", "Designer")
)

$row = 2
foreach ($record in $data) {
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $ws.Cells.Item($row, 4).Value = $record[3]
    $row = $row + 1
}

# Column widths to fit the new text (closest values the character-width
# grid allows to 11.75 / 12.5)
$ws.Columns.Item(3).ColumnWidth = 10.95
$ws.Columns.Item(4).ColumnWidth = 11.8

$ws.Range("A1").Select()
